$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")

# 1. Duplicate Sheet1 to create Sheet2 (placed immediately after Sheet1) while
#    Sheet1 still carries its original values, so every original shared
#    string stays referenced by at least one sheet.
$ws1.Copy([System.Reflection.Missing]::Value, $ws1)
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "Sheet2"

# 2. Update Sheet1 test data (password now matches confirm password).
$ws1.Range("B3").Value = "thanhtung1"
$ws1.Range("C3").Value = "thanhtung1"
$ws1.Range("B4").Value = "ducthuan1"

# Rebuild Sheet1 hyperlinks, reordering D4 ahead of D3
$ws1.Hyperlinks.Delete()
$ws1.Hyperlinks.Add($ws1.Range("D2"), "mailto:viet1@katalon.com")
$ws1.Hyperlinks.Add($ws1.Range("D4"), "mailto:thuanuet@katalon.com")
$ws1.Hyperlinks.Add($ws1.Range("D3"), "mailto:tung1@katalon.com")
$ws1.Range("D2").Style = "Hyperlink"
$ws1.Range("D3").Style = "Hyperlink"
$ws1.Range("D4").Style = "Hyperlink"

$ws1.Range("D9").Select() | Out-Null

# 3. Update Sheet2 test data (new signup test rows).
$ws2.Range("A3").Value = "thanhtung2"
$ws2.Range("A4").Value = "ducthuan2"
$ws2.Range("B4").Value = "ducthuan2"
$ws2.Range("C4").Value = "ducthuan2"
$ws2.Range("D4").ClearContents()

# Rebuild Sheet2 hyperlinks without the D4 entry (D2, D3 only)
$ws2.Hyperlinks.Delete()
$ws2.Hyperlinks.Add($ws2.Range("D2"), "mailto:viet1@katalon.com")
$ws2.Hyperlinks.Add($ws2.Range("D3"), "mailto:tung1@katalon.com")
$ws2.Range("D2").Style = "Hyperlink"
$ws2.Range("D3").Style = "Hyperlink"

$ws2.Activate() | Out-Null
$ws2.Range("A3").Select() | Out-Null
